# Rename the inline picture assets (their display/file "name", i.e. the
# <wp:docPr name="..."> of each embedded image):
#   - both Pearson logo pictures, found in the document's footers:
#       image1.png -> image2.png
#   - the BTec logo picture, found in the document's header:
#       image2.jpg -> image1.jpg
#
# wdHeaderFooterIndex values: 1 = Primary, 2 = FirstPage, 3 = EvenPages

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

for ($i = 1; $i -le 3; $i++) {
    $f = $sec.Footers.Item($i)
    if ($f.Exists) {
        for ($j = 1; $j -le $f.Range.InlineShapes.Count; $j++) {
            $shp = $f.Range.InlineShapes.Item($j)
            $shp.Name = "image2.png"
        }
    }
}

for ($i = 1; $i -le 3; $i++) {
    $h = $sec.Headers.Item($i)
    if ($h.Exists) {
        for ($j = 1; $j -le $h.Range.InlineShapes.Count; $j++) {
            $shp = $h.Range.InlineShapes.Item($j)
            $shp.Name = "image1.jpg"
        }
    }
}
